$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 07:35"

# Hungria (row 72)
$ws.Range("B72").Value = 3713
$ws.Range("C72").Value = 35
$ws.Range("D72").Value = 1655
$ws.Range("E72").Value = 1576
$ws.Range("G72").Value = 6
$ws.Range("H72").Value = 482

# Tailandia (row 76)
$ws.Range("B76").Value = 3040
$ws.Range("C76").Value = 3
$ws.Range("D76").Value = 2916
$ws.Range("E76").Value = 68

# Bulgaria (row 82)
$ws.Range("B82").Value = 2408
$ws.Range("C82").Value = 36
$ws.Range("D82").Value = 808
$ws.Range("E82").Value = 1474
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 126

# Taiwan (row 137)
$ws.Range("D137").Value = 411
$ws.Range("E137").Value = 23
